$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6178.8335
$ws.Range("I32").Value = 3729.849
$ws.Range("K32").Value = 3729.849
$ws.Range("M32").Value = -3442.849
$ws.Range("H61").Value = 3035.5833
$ws.Range("I61").Value = 2200
$ws.Range("J61").Value = 3314.111
$ws.Range("K61").Value = 2200
$ws.Range("L61").Value = 3314.111
$ws.Range("M61").Value = -1988
$ws.Range("N61").Value = -3738.111
$ws.Range("H74").Value = 3310.0625
$ws.Range("I74").Value = 2898.111
$ws.Range("J74").Value = 3839.7144
$ws.Range("K74").Value = 2898.111
$ws.Range("L74").Value = 3839.7144
$ws.Range("M74").Value = -2024.111
$ws.Range("N74").Value = -5587.7144
$ws.Range("H77").Value = 3310.0625
$ws.Range("I77").Value = 2898.111
$ws.Range("J77").Value = 3839.7144
$ws.Range("K77").Value = 14490.555
$ws.Range("L77").Value = 19198.572
$ws.Range("M77").Value = -10122.555
$ws.Range("N77").Value = -27934.572
$ws.Range("H122").Value = 2063.6
$ws.Range("I122").Value = 1201
$ws.Range("K122").Value = 3603
$ws.Range("M122").Value = -1153
$ws.Range("H132").Value = 2423.827
$ws.Range("I132").Value = 1782.5151
$ws.Range("J132").Value = 3537.6843
$ws.Range("K132").Value = 5347.5453
$ws.Range("L132").Value = 10613.0529
$ws.Range("M132").Value = -2817.5453
$ws.Range("N132").Value = -15673.0529
$ws.Range("H136").Value = 3035.5833
$ws.Range("I136").Value = 2200
$ws.Range("J136").Value = 3314.111
$ws.Range("K136").Value = 6600
$ws.Range("L136").Value = 9942.332999999999
$ws.Range("M136").Value = -4050
$ws.Range("N136").Value = -15042.333
$ws.Range("H137").Value = 49674
$ws.Range("J137").Value = 49674
$ws.Range("L137").Value = 49674
$ws.Range("N137").Value = -59874
$ws.Range("H139").Value = 43415.555
$ws.Range("J139").Value = 43415.555
$ws.Range("L139").Value = 43415.555
$ws.Range("N139").Value = -53695.555

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2875982.8
$ws.Range("J7").Value = 25995.4
$ws.Range("L7").Value = 25995.4
$ws.Range("N7").Value = -26221.4
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H94").Value = 776.02563
$ws.Range("I94").Value = 749.0278
$ws.Range("J94").Value = 1100
$ws.Range("K94").Value = 749.0278
$ws.Range("L94").Value = 1100
$ws.Range("M94").Value = -298.0278
$ws.Range("N94").Value = -2002
$ws.Range("H107").Value = 1099.963
$ws.Range("I107").Value = 1022.6818
$ws.Range("J107").Value = 1440
$ws.Range("K107").Value = 1022.6818
$ws.Range("L107").Value = 1440
$ws.Range("M107").Value = 897.3182
$ws.Range("N107").Value = -5280
$ws.Range("H137").Value = 32937.5
$ws.Range("J137").Value = 32937.5
$ws.Range("L137").Value = 32937.5
$ws.Range("N137").Value = -43137.5
$ws.Range("H138").Value = 41334
$ws.Range("J138").Value = 41334
$ws.Range("L138").Value = 41334
$ws.Range("N138").Value = -51614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12823488
$ws.Range("I31").Value = 1205.9565
$ws.Range("J31").Value = 31255518
$ws.Range("K31").Value = 1205.9565
$ws.Range("L31").Value = 31255518
$ws.Range("M31").Value = -910.9565
$ws.Range("N31").Value = -31256108
$ws.Range("H34").Value = 12823488
$ws.Range("I34").Value = 1205.9565
$ws.Range("J34").Value = 31255518
$ws.Range("K34").Value = 1205.9565
$ws.Range("L34").Value = 31255518
$ws.Range("M34").Value = -1003.9565
$ws.Range("N34").Value = -31255922
$ws.Range("H36").Value = 20009.6
$ws.Range("J36").Value = 26666.666
$ws.Range("L36").Value = 26666.666
$ws.Range("N36").Value = -27442.666
$ws.Range("H40").Value = 20009.6
$ws.Range("J40").Value = 26666.666
$ws.Range("L40").Value = 26666.666
$ws.Range("N40").Value = -26986.666
$ws.Range("H42").Value = 42833.2
$ws.Range("J42").Value = 42833.2
$ws.Range("L42").Value = 42833.2
$ws.Range("N42").Value = -44019.2
$ws.Range("H44").Value = 39999.668
$ws.Range("J44").Value = 39999.668
$ws.Range("L44").Value = 39999.668
$ws.Range("N44").Value = -40883.668
$ws.Range("H55").Value = 39999.5
$ws.Range("J55").Value = 39999.5
$ws.Range("L55").Value = 39999.5
$ws.Range("N55").Value = -40629.5
$ws.Range("H138").Value = 29723.334
$ws.Range("J138").Value = 29723.334
$ws.Range("L138").Value = 29723.334
$ws.Range("N138").Value = -40003.334
$ws.Range("H140").Value = 94314.55
$ws.Range("J140").Value = 94314.55
$ws.Range("L140").Value = 94314.55
$ws.Range("N140").Value = -104674.55
$ws.Range("H141").Value = 15476.923
$ws.Range("J141").Value = 15476.923
$ws.Range("L141").Value = 15476.923
$ws.Range("N141").Value = -25836.923

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1508.9615
$ws.Range("I5").Value = 368.35294
$ws.Range("J5").Value = 3663.4443
$ws.Range("K5").Value = 1105.05882
$ws.Range("L5").Value = 10990.3329
$ws.Range("M5").Value = -993.05882
$ws.Range("N5").Value = -11214.3329
$ws.Range("H74").Value = 9669
$ws.Range("I74").Value = 5013
$ws.Range("J74").Value = 11997
$ws.Range("K74").Value = 15039
$ws.Range("L74").Value = 35991
$ws.Range("M74").Value = -13978
$ws.Range("N74").Value = -38113
$ws.Range("H77").Value = 9669
$ws.Range("I77").Value = 5013
$ws.Range("J77").Value = 11997
$ws.Range("K77").Value = 45117
$ws.Range("L77").Value = 107973
$ws.Range("M77").Value = -39813
$ws.Range("N77").Value = -118581
$ws.Range("H82").Value = 6218.8335
$ws.Range("I82").Value = 756.5
$ws.Range("J82").Value = 8950
$ws.Range("K82").Value = 2269.5
$ws.Range("L82").Value = 26850
$ws.Range("M82").Value = -1863.5
$ws.Range("N82").Value = -27662
$ws.Range("H85").Value = 6218.8335
$ws.Range("I85").Value = 756.5
$ws.Range("J85").Value = 8950
$ws.Range("K85").Value = 2269.5
$ws.Range("L85").Value = 26850
$ws.Range("M85").Value = -865.5
$ws.Range("N85").Value = -29658
$ws.Range("H88").Value = 6840
$ws.Range("J88").Value = 6840
$ws.Range("L88").Value = 20520
$ws.Range("N88").Value = -21376
$ws.Range("H91").Value = 6840
$ws.Range("J91").Value = 6840
$ws.Range("L91").Value = 20520
$ws.Range("N91").Value = -23484
$ws.Range("H113").Value = 548.1896400000001
$ws.Range("I113").Value = 544.9286
$ws.Range("J113").Value = 556.75
$ws.Range("K113").Value = 1634.7858
$ws.Range("L113").Value = 1670.25
$ws.Range("M113").Value = 535.2142000000001
$ws.Range("N113").Value = -6010.25
$ws.Range("H123").Value = 2535
$ws.Range("I123").Value = 3202.5
$ws.Range("J123").Value = 1200
$ws.Range("K123").Value = 9607.5
$ws.Range("L123").Value = 3600
$ws.Range("M123").Value = -7157.5
$ws.Range("N123").Value = -8500
$ws.Range("H135").Value = 1508.9615
$ws.Range("I135").Value = 368.35294
$ws.Range("J135").Value = 3663.4443
$ws.Range("K135").Value = 3315.17646
$ws.Range("L135").Value = 32970.9987
$ws.Range("M135").Value = -780.1764599999997
$ws.Range("N135").Value = -38040.9987

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 32304
$ws.Range("J46").Value = 32304
$ws.Range("L46").Value = 32304
$ws.Range("N46").Value = -32616
$ws.Range("H107").Value = 915
$ws.Range("I107").Value = 623
$ws.Range("K107").Value = 623
$ws.Range("M107").Value = 1297
$ws.Range("H137").Value = 63674.6
$ws.Range("J137").Value = 63674.6
$ws.Range("L137").Value = 63674.6
$ws.Range("N137").Value = -73874.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7770.5713
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 8565.666999999999
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 8565.666999999999
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -8789.666999999999
$ws.Range("H40").Value = 8237.15
$ws.Range("I40").Value = 7458.154
$ws.Range("K40").Value = 7458.154
$ws.Range("M40").Value = -7322.154
$ws.Range("H46").Value = 1923.28
$ws.Range("I46").Value = 1637.625
$ws.Range("J46").Value = 2431.111
$ws.Range("K46").Value = 1637.625
$ws.Range("L46").Value = 2431.111
$ws.Range("M46").Value = -1449.625
$ws.Range("N46").Value = -2807.111
$ws.Range("H126").Value = 7770.5713
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 8565.666999999999
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 25697.001
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -30637.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31248
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96240
$ws.Range("H132").Value = 14495157
$ws.Range("I132").Value = 756.4
$ws.Range("K132").Value = 2269.2
$ws.Range("M132").Value = 260.8000000000002
